# Update 2021 HWL2 First Batch
# - Extend the year range from 2015 through 2050 on the "Data Clio Infra
#   Format" and "Data Long Format" sheets.
# - Fix the Zwart/van Leeuwen citation text on the Metadata sheet (added a
#   comma before "and Jieli van Leeuwen-Li").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Data Clio Infra Format": the sheet already has year columns running
#    from 1500 (column I) to 2015 (column TD). Append 35 more year columns,
#    TE:UM, for 2016..2050, as header text in row 1 only (row 2's new cells
#    stay blank, matching the source data which has no values past 2015).
# ---------------------------------------------------------------------
$wsClio = $wb.Worksheets.Item("Data Clio Infra Format")

$clioFirstCol = 525   # column TE
$clioHeaderRange = $wsClio.Range("TE1:UM1")
$clioHeaderRange.NumberFormat = "@"
for ($i = 0; $i -lt 35; $i++) {
    $col = $clioFirstCol + $i
    $year = 2016 + $i
    $wsClio.Cells.Item(1, $col).Value = "$year"
}
$clioHeaderRange.ClearFormats()

# ---------------------------------------------------------------------
# 2) "Data Long Format": insert 35 new columns right before the existing
#    "year"/"value" columns (currently E:F) and give them the same
#    2016..2050 headers. This pushes the old E/F ("year"/"value") columns
#    out to AN/AO, carrying their data along automatically.
# ---------------------------------------------------------------------
$wsLong = $wb.Worksheets.Item("Data Long Format")

$wsLong.Range("E1:AM1").EntireColumn.Insert()

$longFirstCol = 5   # column E
$longHeaderRange = $wsLong.Range("E1:AM1")
$longHeaderRange.NumberFormat = "@"
for ($i = 0; $i -lt 35; $i++) {
    $col = $longFirstCol + $i
    $year = 2016 + $i
    $wsLong.Cells.Item(1, $col).Value = "$year"
}
$longHeaderRange.ClearFormats()

# ---------------------------------------------------------------------
# 3) "Metadata": fix the citation text (added a comma after "van Leeuwen").
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("C3").Value = "Zwart, Pim de, Bas van Leeuwen, and Jieli van Leeuwen-Li (2015). Labourers Real Wage. http://hdl.handle.net/10622/QK8VRF, accessed via the Clio Infra website."

Write-Host "Applied 2016-2050 year columns and citation fix"
